$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function Set-ThemeColor($index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $rgb = $r + ($g * 256) + ($b * 65536)
    $tcs.Colors($index).RGB = $rgb
}

Set-ThemeColor 1  "000000"
Set-ThemeColor 2  "FFFFFF"
Set-ThemeColor 3  "44546A"
Set-ThemeColor 4  "E7E6E6"
Set-ThemeColor 5  "5B9BD5"
Set-ThemeColor 6  "ED7D31"
Set-ThemeColor 7  "A5A5A5"
Set-ThemeColor 8  "FFC000"
Set-ThemeColor 9  "4472C4"
Set-ThemeColor 10 "70AD47"
Set-ThemeColor 11 "0563C1"
Set-ThemeColor 12 "954F72"
